# Add a new worksheet named "rho" as the last sheet in the workbook,
# matching the structure/style used by the existing "Adjusted2" sheet,
# and fill it with the theta/Jxx/Jyy/beta/gamma/trace_sq table.

$wb = $excel.ActiveWorkbook

# Use the last existing worksheet as a style template (same header
# formatting: bold, centered, bordered) and as the "after" anchor so the
# new sheet is appended at the end of the tab order.
$templateSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $templateSheet)
$newSheet.Name = "rho"

# Copy header formatting (bold font, centered alignment, thin border)
# from the template sheet's header row.
$templateSheet.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row
$newSheet.Range("A1").Value = "theta"
$newSheet.Range("B1").Value = "Jxx"
$newSheet.Range("C1").Value = "Jyy"
$newSheet.Range("D1").Value = "beta"
$newSheet.Range("E1").Value = "gamma"
$newSheet.Range("F1").Value = "trace_sq"

$data = @(
  @(0,0.000506817080977578,0.9994933377741827,-0.004801290864711461,0,0.9990332939064652),
  @(5,0.03154104479697824,0.9684590087845358,-0.1717602611156163,0,0.9979108637998201),
  @(10,0.1214340627213908,0.8785664893102039,-0.3217326011874075,0,0.99364904106151),
  @(15,0.2517583769675233,0.7482419850665917,-0.4269370214922129,0,0.9877987892309996),
  @(20,0.424158175943422,0.5758418789718585,-0.4855717765803552,0,0.9830639282202962),
  @(25,0.5939355094920483,0.4060644546166343,-0.4748352469848616,0,0.9685847542970327),
  @(30,0.7437731041825903,0.2562254010684298,-0.4160468140586621,0,0.9650397896348101),
  @(35,0.8718868327665467,0.1281131623507801,-0.3046338543527607,0,0.9622032019548354),
  @(40,0.950662521476947,0.04933725337367501,-0.1554673430906546,0,0.9545335838466994),
  @(45,0.9778614205294466,0.02213860108200471,0.01034475739352828,0,0.9569171034287971),
  @(50,0.9454666834656507,0.05453331352093595,0.1743684209333196,0,0.9576898242646681),
  @(55,0.8668700253609699,0.1331299456045153,0.3131342195039139,0,0.9652933021346404),
  @(60,0.7362297773998635,0.2637702257249902,0.4238350725941739,0,0.9708813546310823),
  @(65,0.5808107695672143,0.4191892603515385,0.4819689331889778,0,0.9776488911579724),
  @(70,0.4154905655747549,0.5845093624310944,0.4831127830167778,0,0.9810795270796667),
  @(75,0.2548624967861255,0.7451378698426443,0.4287498866127726,0,0.9878382678826222),
  @(80,0.1187095330880051,0.8812910806824168,0.3185539706103371,0,0.9937191865195772),
  @(85,0.03333951007432865,0.9666604933782167,0.1758887008222358,0,0.9974177025440816),
  @(90,0.0007035867102550072,0.9992963996430539,-0.002401609023337846,0,0.998605324825631)
)

$r = 2
foreach ($row in $data) {
  $newSheet.Cells.Item($r, 1).Value = $row[0]
  $newSheet.Cells.Item($r, 2).Value = $row[1]
  $newSheet.Cells.Item($r, 3).Value = $row[2]
  $newSheet.Cells.Item($r, 4).Value = $row[3]
  $newSheet.Cells.Item($r, 5).Value = $row[4]
  $newSheet.Cells.Item($r, 6).Value = $row[5]
  $r++
}

$newSheet.Range("A1").Select()
